$wb = $excel.ActiveWorkbook

$wsGen  = $wb.Worksheets.Item("Generators")
$wsLoad = $wb.Worksheets.Item("Load")

# 1) Fix the xf currently sitting at cellXfs index 5 (owned by B2) IN PLACE:
#    horizontal-center -> vertical-center.
$wsLoad.Range("B2").HorizontalAlignment = 1      # xlGeneral
$wsLoad.Range("B2").VerticalAlignment = -4108    # xlCenter

# 2) Give C1 the same centered header style as B1 (cellXfs index 4).
$wsLoad.Range("B1").Copy()
$wsLoad.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$wsLoad.Application.CutCopyMode = $false

# 2b) A1 switches from the centered header style to the plain (non-centered) header
#     style used on the Generators sheet header row (cellXfs index 1).
$wsGen.Range("A1").Copy()
$wsLoad.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$wsLoad.Application.CutCopyMode = $false

# 3) Give A2:A3 the numeric centered style (cellXfs index 3), matching Generators!A2.
$wsGen.Range("A2").Copy()
$wsLoad.Range("A2:A3").PasteSpecial(-4122)
$wsLoad.Application.CutCopyMode = $false

# 4) Move the "500"-style formatting that lived on B2 over to the new C2:C3 column.
$wsLoad.Range("B2").Copy()
$wsLoad.Range("C2:C3").PasteSpecial(-4122)
$wsLoad.Application.CutCopyMode = $false

# 5) Column B becomes the plain "Time" column now - clear its leftover formatting.
$wsLoad.Range("B2:B3").ClearFormats()

# 6) Write the final cell values for the extended "Load" table (Id, Time, Capacity).
$wsLoad.Range("A1").Value = "Id"
$wsLoad.Range("B1").Value = "Time"
$wsLoad.Range("C1").Value = "Capacity"

$wsLoad.Range("A2").Value = 1
$wsLoad.Range("B2").Value = 1
$wsLoad.Range("C2").Value = 500

$wsLoad.Range("A3").Value = 2
$wsLoad.Range("B3").Value = 2
$wsLoad.Range("C3").Value = 450

# 7) Selections: Generators gets a new A1:A3 selection (it stays the inactive tab),
#    Load (the active tab) moves its selection to H6.
$wsGen.Activate()
$wsGen.Range("A1:A3").Select() | Out-Null

$wsLoad.Activate()
$wsLoad.Range("H6").Select() | Out-Null

# 8) Page setup on the Load sheet: portrait orientation.
$wsLoad.PageSetup.Orientation = 1   # xlPortrait
